# Junction_Flooding_233.xlsx edit:
#  - Round row 5 values down to a "custom accuracy" of 2 decimal places
#  - Remove row 6 (reduces the data set / dimension to A1:AH5)
#  - A handful of columns end up one character narrower as a side effect
#    of the values becoming shorter after rounding

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to their 2-decimal rounded equivalents ---
$ws.Range("B5").Value = 18.15
$ws.Range("C5").Value = 13.48
$ws.Range("D5").Value = 1.16
$ws.Range("E5").Value = 39.71
$ws.Range("F5").Value = 32.21
$ws.Range("G5").Value = 14.23
$ws.Range("H5").Value = 53.71
$ws.Range("I5").Value = 22.11
$ws.Range("J5").Value = 9.76
$ws.Range("K5").Value = 14.38
$ws.Range("L5").Value = 15.92
$ws.Range("M5").Value = 16.95
$ws.Range("N5").Value = 4.59
$ws.Range("O5").Value = 14.29
$ws.Range("P5").Value = 20.27
$ws.Range("Q5").Value = 12.13
$ws.Range("R5").Value = 0.77
$ws.Range("S5").Value = 0.77
$ws.Range("T5").Value = 210.13
$ws.Range("U5").Value = 39.9
$ws.Range("V5").Value = 13.19
$ws.Range("W5").Value = 26.73
$ws.Range("X5").Value = 14.02
$ws.Range("Y5").Value = 2.17
$ws.Range("Z5").Value = 26.44
$ws.Range("AA5").Value = 11.65
$ws.Range("AB5").Value = 10.36
$ws.Range("AC5").Value = 12.19
$ws.Range("AD5").Value = 16.7
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 48.81
$ws.Range("AG5").Value = 7.38
$ws.Range("AH5").Value = 16.49

# --- Remove row 6 entirely (shrinks used range to A1:AH5) ---
$ws.Rows.Item(6).Delete()

# --- Narrow the columns whose best-fit width shrank after the rounding ---
$ws.Columns.Item(3).ColumnWidth = 6.15
$ws.Columns.Item(7).ColumnWidth = 6.15
$ws.Columns.Item(11).ColumnWidth = 6.15
$ws.Columns.Item(17).ColumnWidth = 6.15
$ws.Columns.Item(22).ColumnWidth = 6.15
$ws.Columns.Item(27).ColumnWidth = 6.15
$ws.Columns.Item(28).ColumnWidth = 6.15
$ws.Columns.Item(29).ColumnWidth = 6.15
